# Change the deck's theme from the custom "Integral" (Red Violet) colour
# scheme to the stock PowerPoint "Office Theme" (Office) colour scheme —
# i.e. what happens when Design > Themes > Office Theme is picked in the
# ribbon. Colours are written through the Theme's ThemeColorScheme, which
# is the supported, persisted way to edit a theme's 12 colour slots via
# the PowerPoint object model (RGB is stored OLE-style as 0xBBGGRR).

function New-Bgr($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# The 12 colour slots of the "Office" colour scheme, in
# ThemeColorScheme.Item() order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    (New-Bgr 0x00 0x00 0x00),   # dk1      000000
    (New-Bgr 0xFF 0xFF 0xFF),   # lt1      FFFFFF
    (New-Bgr 0x44 0x54 0x6A),   # dk2      44546A
    (New-Bgr 0xE7 0xE6 0xE6),   # lt2      E7E6E6
    (New-Bgr 0x5B 0x9B 0xD5),   # accent1  5B9BD5
    (New-Bgr 0xED 0x7D 0x31),   # accent2  ED7D31
    (New-Bgr 0xA5 0xA5 0xA5),   # accent3  A5A5A5
    (New-Bgr 0xFF 0xC0 0x00),   # accent4  FFC000
    (New-Bgr 0x44 0x72 0xC4),   # accent5  4472C4
    (New-Bgr 0x70 0xAD 0x47),   # accent6  70AD47
    (New-Bgr 0x05 0x63 0xC1),   # hlink    0563C1
    (New-Bgr 0x95 0x4F 0x72)    # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}

Write-Output "Applied Office Theme colour scheme to the presentation theme."
